$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-18 Wednesday" "2025-06-19 Thursday"

Replace-Text "139÷6=" "719÷7="
Replace-Text "538÷4=" "616÷8="
Replace-Text "200÷3=" "919÷2="
Replace-Text "117÷8=" "820÷5="
Replace-Text "477÷3=" "232÷4="
Replace-Text "443÷7=" "395÷9="
Replace-Text "282÷9=" "251÷2="
Replace-Text "958÷8=" "304÷2="
Replace-Text "564÷6=" "708÷4="
Replace-Text "860÷3=" "185÷9="
Replace-Text "288÷6=" "678÷6="
Replace-Text "193÷6=" "196÷2="
Replace-Text "258÷7=" "457÷2="
Replace-Text "814÷9=" "654÷4="
Replace-Text "470÷3=" "187÷2="
Replace-Text "388÷2=" "269÷5="
Replace-Text "113÷4=" "656÷5="
Replace-Text "397÷6=" "468÷7="
Replace-Text "675÷4=" "792÷9="
Replace-Text "905÷4=" "308÷9="
Replace-Text "937÷6=" "867÷6="
Replace-Text "216÷4=" "348÷2="
Replace-Text "978÷5=" "200÷7="
Replace-Text "805÷8=" "574÷9="
Replace-Text "910÷8=" "610÷3="
